$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(12978, 12648, 12539, 12135, 12005, 11277, 11277, 11277, 11277, 10859, 10859, 10370, 10074, 10074, 10074, 9039, 9039, 9039, 8971, 8971, 8971, 8971, 8880, 8880, 8880, 8880, 8880, 8880, 8880, 8880, 8880, 8880, 8880, 8880, 8880, 8880, 8822, 8822, 8822, 8292, 8292, 8292, 8292, 8292, 8292, 8292, 8292, 8292, 8292, 8292, 8292, 8292, 8292, 8292, 8292, 8292, 8292, 8292, 8292, 8292, 8292, 8292, 8255, 7787, 7787, 7787, 7787, 7787, 7787, 7787, 7787, 7787, 7787, 7787, 7787, 7787, 7787, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7295, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $values[$i]
}
